$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03427420202738
$ws.Range("D2").Value = 1.041365459208463
$ws.Range("E2").Value = 1.037918556866448
$ws.Range("F2").Value = 1.048674377599796
$ws.Range("I2").Value = 1.033796666875025
$ws.Range("J2").Value = 1.039394344286892
$ws.Range("K2").Value = 1.04414517993126
$ws.Range("L2").Value = 1.040708070665837
$ws.Range("M2").Value = 1.051433563782909
$ws.Range("N2").Value = 1.040870403475977

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035413281770566
$ws.Range("D3").Value = 1.042257216117441
$ws.Range("E3").Value = 1.039006164679349
$ws.Range("F3").Value = 1.049765546586225
$ws.Range("I3").Value = 1.033997551711782
$ws.Range("J3").Value = 1.040175527966717
$ws.Range("K3").Value = 1.044847336605564
$ws.Range("L3").Value = 1.041604850660635
$ws.Range("M3").Value = 1.052336103849365
$ws.Range("N3").Value = 1.041652696526231

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036150417961636
$ws.Range("D4").Value = 1.042834079777768
$ws.Range("E4").Value = 1.039710354654507
$ws.Range("F4").Value = 1.050471906939601
$ws.Range("I4").Value = 1.034126084261762
$ws.Range("J4").Value = 1.040680554787859
$ws.Range("K4").Value = 1.045300883760595
$ws.Range("L4").Value = 1.042184972316799
$ws.Range("M4").Value = 1.052919816273814
$ws.Range("N4").Value = 1.042158440543385

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.036460329013563
$ws.Range("D5").Value = 1.04307655424501
$ws.Range("E5").Value = 1.040006500622412
$ws.Range("F5").Value = 1.050768933130318
$ws.Range("I5").Value = 1.034179771317323
$ws.Range("J5").Value = 1.040892760305759
$ws.Range("K5").Value = 1.045491364683521
$ws.Range("L5").Value = 1.042428818672646
$ws.Range("M5").Value = 1.053165139402035
$ws.Range("N5").Value = 1.042370947417462

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036512365584242
$ws.Range("D6").Value = 1.043117264485134
$ws.Range("E6").Value = 1.040056230972208
$ws.Range("F6").Value = 1.050818809380879
$ws.Range("I6").Value = 1.034188765199069
$ws.Range("J6").Value = 1.040928384234975
$ws.Range("K6").Value = 1.04552333611112
$ws.Range("L6").Value = 1.042469759400537
$ws.Range("M6").Value = 1.053206326159912
$ws.Range("N6").Value = 1.042406621936743

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036154558932451
$ws.Range("D7").Value = 1.042837319886862
$ws.Range("E7").Value = 1.039714311360345
$ws.Range("F7").Value = 1.050475875533735
$ws.Range("I7").Value = 1.03412680299845
$ws.Range("J7").Value = 1.040683390711449
$ws.Range("K7").Value = 1.04530342972369
$ws.Range("L7").Value = 1.042188230747906
$ws.Range("M7").Value = 1.052923094566627
$ws.Range("N7").Value = 1.042161280494312

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034659144306408
$ws.Range("D8").Value = 1.041666866277503
$ws.Range("E8").Value = 1.038286029169804
$ws.Range("F8").Value = 1.049043081068038
$ws.Range("I8").Value = 1.033864857723645
$ws.Range("J8").Value = 1.03965844256966
$ws.Range("K8").Value = 1.044382641749613
$ws.Range("L8").Value = 1.041011173937972
$ws.Range("M8").Value = 1.051738641855523
$ws.Range("N8").Value = 1.041134876808599

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.032024568046465
$ws.Range("D9").Value = 1.039603133810593
$ws.Range("E9").Value = 1.035772528221528
$ws.Range("F9").Value = 1.046520609235137
$ws.Range("I9").Value = 1.033392149250579
$ws.Range("J9").Value = 1.037848884211532
$ws.Range("K9").Value = 1.042754002329342
$ws.Range("L9").Value = 1.038935848220008
$ws.Range("M9").Value = 1.049649245096179
$ws.Range("N9").Value = 1.039322748670044

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030268482115438
$ws.Range("D10").Value = 1.038226476145424
$ws.Range("E10").Value = 1.034099063192146
$ws.Range("F10").Value = 1.04484048637901
$ws.Range("I10").Value = 1.033069532583194
$ws.Range("J10").Value = 1.036640159565608
$ws.Range("K10").Value = 1.04166413859853
$ws.Range("L10").Value = 1.037551467366822
$ws.Range("M10").Value = 1.048254793471264
$ws.Range("N10").Value = 1.038112307496479

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029508133721428
$ws.Range("D11").Value = 1.037630168934448
$ws.Range("E11").Value = 1.03337494991754
$ws.Range("F11").Value = 1.044113331340214
$ws.Range("I11").Value = 1.032928062795348
$ws.Range("J11").Value = 1.036116205222174
$ws.Range("K11").Value = 1.041191239683976
$ws.Range("L11").Value = 1.036951812328442
$ws.Range("M11").Value = 1.047650615739737
$ws.Range("N11").Value = 1.03758760907778

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029225712391968
$ws.Range("D12").Value = 1.037408642838897
$ws.Range("E12").Value = 1.033106057657404
$ws.Range("F12").Value = 1.043843285508089
$ws.Range("I12").Value = 1.032875247966356
$ws.Range("J12").Value = 1.035921499275692
$ws.Range("K12").Value = 1.041015436314978
$ws.Range("L12").Value = 1.036729041741522
$ws.Range("M12").Value = 1.0474261410282
$ws.Range("N12").Value = 1.03739262662652

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029286292479737
$ws.Range("D13").Value = 1.037456162346462
$ws.Range("E13").Value = 1.033163732547681
$ws.Range("F13").Value = 1.043901208914774
$ws.Range("I13").Value = 1.03288658899347
$ws.Range("J13").Value = 1.035963268261505
$ws.Range("K13").Value = 1.041053153421212
$ws.Range("L13").Value = 1.03677682824673
$ws.Range("M13").Value = 1.047474294166537
$ws.Range("N13").Value = 1.037434454929084

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029484788573318
$ws.Range("D14").Value = 1.037611858152762
$ws.Range("E14").Value = 1.033352721660161
$ws.Range("F14").Value = 1.044091008211369
$ws.Range("I14").Value = 1.032923702544883
$ws.Range("J14").Value = 1.036100112524927
$ws.Range("K14").Value = 1.041176710736612
$ws.Range("L14").Value = 1.03693339869866
$ws.Range("M14").Value = 1.04763206173924
$ws.Range("N14").Value = 1.037571493527056

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02960708931158
$ws.Range("D15").Value = 1.037707783357778
$ws.Range("E15").Value = 1.033469174093683
$ws.Range("F15").Value = 1.04420795668655
$ws.Range("I15").Value = 1.03294653409406
$ws.Range("J15").Value = 1.036184415383186
$ws.Range("K15").Value = 1.041252818901252
$ws.Range("L15").Value = 1.03702986266229
$ws.Range("M15").Value = 1.047729260086625
$ws.Range("N15").Value = 1.037655916105044

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030318944838451
$ws.Range("D16").Value = 1.038266046747966
$ws.Range("E16").Value = 1.034147130852217
$ws.Range("F16").Value = 1.044888752575873
$ws.Range("I16").Value = 1.033078884074789
$ws.Range("J16").Value = 1.036674920663268
$ws.Range("K16").Value = 1.041695502652547
$ws.Range("L16").Value = 1.03759126004175
$ws.Range("M16").Value = 1.048294882901702
$ws.Range("N16").Value = 1.038147117958885

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030765485007056
$ws.Range("D17").Value = 1.038616175507524
$ws.Range("E17").Value = 1.034572530769544
$ws.Range("F17").Value = 1.045315890970228
$ws.Range("I17").Value = 1.03316142857732
$ws.Range("J17").Value = 1.036982449050226
$ws.Range("K17").Value = 1.04197292368538
$ws.Range("L17").Value = 1.037943353404506
$ws.Range("M17").Value = 1.048649583500945
$ws.Range("N17").Value = 1.038455083071427

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031025949093436
$ws.Range("D18").Value = 1.038820379888436
$ws.Range("E18").Value = 1.034820708417233
$ws.Range("F18").Value = 1.045565067302724
$ws.Range("I18").Value = 1.033209404163975
$ws.Range("J18").Value = 1.037161770331058
$ws.Range("K18").Value = 1.042134644075384
$ws.Range("L18").Value = 1.038148703382001
$ws.Range("M18").Value = 1.048856438596966
$ws.Range("N18").Value = 1.038634659009047

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031114761464602
$ws.Range("D19").Value = 1.038890004952515
$ws.Range("E19").Value = 1.034905338926203
$ws.Range("F19").Value = 1.045650035731162
$ws.Range("I19").Value = 1.033225733557413
$ws.Range("J19").Value = 1.037222904971284
$ws.Range("K19").Value = 1.042189770512027
$ws.Range("L19").Value = 1.038218719022673
$ws.Range("M19").Value = 1.04892696478732
$ws.Range("N19").Value = 1.038695880467473

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030717574990099
$ws.Range("D20").Value = 1.038578612054565
$ws.Range("E20").Value = 1.034526884313302
$ws.Range("F20").Value = 1.04527005957744
$ws.Range("I20").Value = 1.03315259004418
$ws.Range("J20").Value = 1.036949459839285
$ws.Range("K20").Value = 1.041943168821599
$ws.Range("L20").Value = 1.037905579204616
$ws.Range("M20").Value = 1.048611531185791
$ws.Range("N20").Value = 1.038422047012022

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029426336268627
$ws.Range("D21").Value = 1.037566010477039
$ws.Range("E21").Value = 1.03329706699416
$ws.Range("F21").Value = 1.044035115620114
$ws.Range("I21").Value = 1.032912780890141
$ws.Range("J21").Value = 1.036059817677829
$ws.Range("K21").Value = 1.041140330265371
$ws.Range("L21").Value = 1.036887293491386
$ws.Range("M21").Value = 1.047585604676505
$ws.Range("N21").Value = 1.037531141456653

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028614516809053
$ws.Range("D22").Value = 1.036929167949898
$ws.Range("E22").Value = 1.032524269387039
$ws.Range("F22").Value = 1.043258957496014
$ws.Range("I22").Value = 1.032760460582057
$ws.Range("J22").Value = 1.035499966648004
$ws.Range("K22").Value = 1.040634699350537
$ws.Range("L22").Value = 1.036246871369665
$ws.Range("M22").Value = 1.046940238567699
$ws.Range("N22").Value = 1.03697049537415

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029044874934163
$ws.Range("D23").Value = 1.0372667873245
$ws.Range("E23").Value = 1.032933902794849
$ws.Range("F23").Value = 1.043670385333814
$ws.Range("I23").Value = 1.03284135470701
$ws.Range("J23").Value = 1.035796801630948
$ws.Range("K23").Value = 1.040902824962385
$ws.Range("L23").Value = 1.036586389037839
$ws.Range("M23").Value = 1.047282390352338
$ws.Range("N23").Value = 1.037267751896818

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030739223449379
$ws.Range("D24").Value = 1.038595585423057
$ws.Range("E24").Value = 1.034547509830143
$ws.Range("F24").Value = 1.045290768706562
$ws.Range("I24").Value = 1.033156584325841
$ws.Range("J24").Value = 1.036964366412869
$ws.Range("K24").Value = 1.041956614055881
$ws.Range("L24").Value = 1.037922647801502
$ws.Range("M24").Value = 1.048628725499332
$ws.Range("N24").Value = 1.038436974754651

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03270561207316
$ws.Range("D25").Value = 1.04013680547366
$ws.Range("E25").Value = 1.036421938361592
$ws.Range("F25").Value = 1.047172458540555
$ws.Range("I25").Value = 1.033515673988967
$ws.Range("J25").Value = 1.038317111564685
$ws.Range("K25").Value = 1.043175767494788
$ws.Range("L25").Value = 1.039472514351653
$ws.Range("M25").Value = 1.050189671156413
$ws.Range("N25").Value = 1.039791640959745
